$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 458
$ws.Range("F5").Value = 332
$ws.Range("F6").Value = 476
$ws.Range("F8").Value = 2214
$ws.Range("F9").Value = 53
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 1657
$ws.Range("F12").Value = 1657
$ws.Range("F13").Value = 1370
$ws.Range("F15").Value = 1424
$ws.Range("F18").Value = 596
$ws.Range("F19").Value = 164
$ws.Range("F21").Value = 7305
$ws.Range("F22").Value = 8104
$ws.Range("F23").Value = 51
$ws.Range("F25").Value = 205
$ws.Range("F34").Value = 350
$ws.Range("F35").Value = 1452
$ws.Range("F36").Value = 247
$ws.Range("F38").Value = 15
$ws.Range("F40").Value = 16
$ws.Range("F41").Value = 744
$ws.Range("F45").Value = 247
$ws.Range("F46").Value = 206
$ws.Range("F48").Value = 182
$ws.Range("F49").Value = 166
$ws.Range("F50").Value = 16

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 8
$ws.Range("F5").Value = 62
$ws.Range("F18").Value = 300

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2626
$ws.Range("F4").Value = 285
$ws.Range("F6").Value = 16

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 458
$ws.Range("F6").Value = 285
$ws.Range("F8").Value = 332
$ws.Range("F10").Value = 476
$ws.Range("F11").Value = 53
$ws.Range("F12").Value = 64
$ws.Range("F13").Value = 1657
$ws.Range("F14").Value = 1657
$ws.Range("F15").Value = 16
$ws.Range("F16").Value = 1370
$ws.Range("F20").Value = 596
$ws.Range("F21").Value = 164
$ws.Range("F24").Value = 7305
$ws.Range("F25").Value = 8104
$ws.Range("F26").Value = 51
$ws.Range("F27").Value = 205
$ws.Range("F31").Value = 350
$ws.Range("F32").Value = 1452
$ws.Range("F33").Value = 247
$ws.Range("F35").Value = 15
$ws.Range("F39").Value = 744
$ws.Range("F45").Value = 247
$ws.Range("F46").Value = 206
$ws.Range("F47").Value = 182
$ws.Range("F48").Value = 166
$ws.Range("F50").Value = 300
